$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2899.8  # H40: 3071.2856 -> 2899.8
$ws.Cells.Item(40, 10).Value = 0  # J40: 3500 -> 0
$ws.Cells.Item(40, 12).Value = 0  # L40: 3500 -> 0
$ws.Cells.Item(40, 14).ClearContents()  # N40: -3850 -> (removed)
$ws.Cells.Item(58, 8).Value = 15402.667  # H58: 17583.334 -> 15402.667
$ws.Cells.Item(58, 9).Value = 915  # I58: 0 -> 915
$ws.Cells.Item(58, 10).Value = 18300.2  # J58: 17583.334 -> 18300.2
$ws.Cells.Item(58, 11).Value = 2745  # K58: 0 -> 2745
$ws.Cells.Item(58, 12).Value = 54900.60000000001  # L58: 52750.00199999999 -> 54900.60000000001
$ws.Cells.Item(58, 13).Value = -2595  # M58: None -> -2595
$ws.Cells.Item(58, 14).Value = -55200.60000000001  # N58: -53050.00199999999 -> -55200.60000000001
$ws.Cells.Item(70, 8).Value = 1615.3334  # H70: 1600.375 -> 1615.3334
$ws.Cells.Item(70, 9).Value = 1531.5714  # I70: 1497.8889 -> 1531.5714
$ws.Cells.Item(70, 10).Value = 1688.625  # J70: 1732.1428 -> 1688.625
$ws.Cells.Item(70, 11).Value = 4594.7142  # K70: 4493.6667 -> 4594.7142
$ws.Cells.Item(70, 12).Value = 5065.875  # L70: 5196.428400000001 -> 5065.875
$ws.Cells.Item(70, 13).Value = -4324.7142  # M70: -4223.6667 -> -4324.7142
$ws.Cells.Item(70, 14).Value = -5605.875  # N70: -5736.428400000001 -> -5605.875
$ws.Cells.Item(73, 8).Value = 1615.3334  # H73: 1600.375 -> 1615.3334
$ws.Cells.Item(73, 9).Value = 1531.5714  # I73: 1497.8889 -> 1531.5714
$ws.Cells.Item(73, 10).Value = 1688.625  # J73: 1732.1428 -> 1688.625
$ws.Cells.Item(73, 11).Value = 4594.7142  # K73: 4493.6667 -> 4594.7142
$ws.Cells.Item(73, 12).Value = 5065.875  # L73: 5196.428400000001 -> 5065.875
$ws.Cells.Item(73, 13).Value = -3658.7142  # M73: -3557.6667 -> -3658.7142
$ws.Cells.Item(73, 14).Value = -6937.875  # N73: -7068.428400000001 -> -6937.875
$ws.Cells.Item(74, 8).Value = 3899.5  # H74: 4279.6 -> 3899.5
$ws.Cells.Item(74, 9).Value = 2599.25  # I74: 2799.3333 -> 2599.25
$ws.Cells.Item(74, 11).Value = 2599.25  # K74: 2799.3333 -> 2599.25
$ws.Cells.Item(74, 13).Value = -1663.25  # M74: -1863.3333 -> -1663.25
$ws.Cells.Item(76, 8).Value = 6568.6  # H76: 6569.1 -> 6568.6
$ws.Cells.Item(76, 9).Value = 5390  # I76: 0 -> 5390
$ws.Cells.Item(76, 10).Value = 6699.5557  # J76: 6569.1 -> 6699.5557
$ws.Cells.Item(76, 11).Value = 5390  # K76: 0 -> 5390
$ws.Cells.Item(76, 12).Value = 6699.5557  # L76: 6569.1 -> 6699.5557
$ws.Cells.Item(76, 13).Value = -5075  # M76: None -> -5075
$ws.Cells.Item(76, 14).Value = -7329.5557  # N76: -7199.1 -> -7329.5557
$ws.Cells.Item(77, 8).Value = 3899.5  # H77: 4279.6 -> 3899.5
$ws.Cells.Item(77, 9).Value = 2599.25  # I77: 2799.3333 -> 2599.25
$ws.Cells.Item(77, 11).Value = 12996.25  # K77: 13996.6665 -> 12996.25
$ws.Cells.Item(77, 13).Value = -8316.25  # M77: -9316.666499999999 -> -8316.25
$ws.Cells.Item(79, 8).Value = 6568.6  # H79: 6569.1 -> 6568.6
$ws.Cells.Item(79, 9).Value = 5390  # I79: 0 -> 5390
$ws.Cells.Item(79, 10).Value = 6699.5557  # J79: 6569.1 -> 6699.5557
$ws.Cells.Item(79, 11).Value = 5390  # K79: 0 -> 5390
$ws.Cells.Item(79, 12).Value = 6699.5557  # L79: 6569.1 -> 6699.5557
$ws.Cells.Item(79, 13).Value = -4298  # M79: None -> -4298
$ws.Cells.Item(79, 14).Value = -8883.555700000001  # N79: -8753.1 -> -8883.555700000001
$ws.Cells.Item(98, 8).Value = 2638.4119  # H98: 2740.9375 -> 2638.4119
$ws.Cells.Item(98, 9).Value = 2709.5625  # I98: 2823.6667 -> 2709.5625
$ws.Cells.Item(98, 11).Value = 2709.5625  # K98: 2823.6667 -> 2709.5625
$ws.Cells.Item(98, 13).Value = -1211.5625  # M98: -1325.6667 -> -1211.5625
$ws.Cells.Item(116, 8).Value = 3377.516  # H116: 3349.606 -> 3377.516
$ws.Cells.Item(116, 9).Value = 3379.2964  # I116: 3353.2144 -> 3379.2964
$ws.Cells.Item(116, 10).Value = 3365.5  # J116: 3329.4 -> 3365.5
$ws.Cells.Item(116, 11).Value = 3379.2964  # K116: 3353.2144 -> 3379.2964
$ws.Cells.Item(116, 12).Value = 3365.5  # L116: 3329.4 -> 3365.5
$ws.Cells.Item(116, 13).Value = 62.70359999999982  # M116: 88.78560000000016 -> 62.70359999999982
$ws.Cells.Item(116, 14).Value = -10249.5  # N116: -10213.4 -> -10249.5
$ws.Cells.Item(122, 8).Value = 2638.4119  # H122: 2740.9375 -> 2638.4119
$ws.Cells.Item(122, 9).Value = 2709.5625  # I122: 2823.6667 -> 2709.5625
$ws.Cells.Item(122, 11).Value = 8128.6875  # K122: 8471.000100000001 -> 8128.6875
$ws.Cells.Item(122, 13).Value = -5678.6875  # M122: -6021.000100000001 -> -5678.6875
$ws.Cells.Item(132, 8).Value = 18869078  # H132: 17858256 -> 18869078
$ws.Cells.Item(132, 9).Value = 18869078  # I132: 17858256 -> 18869078
$ws.Cells.Item(132, 11).Value = 56607234  # K132: 53574768 -> 56607234
$ws.Cells.Item(132, 13).Value = -56604704  # M132: -53572238 -> -56604704
$ws.Cells.Item(137, 8).Value = 3385.4666  # H137: 3756.4614 -> 3385.4666
$ws.Cells.Item(137, 9).Value = 1955.091  # I137: 2173.111 -> 1955.091
$ws.Cells.Item(137, 11).Value = 5865.272999999999  # K137: 6519.333 -> 5865.272999999999
$ws.Cells.Item(137, 13).Value = -3315.272999999999  # M137: -3969.333 -> -3315.272999999999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2301.2  # H2: 2347.2 -> 2301.2
$ws.Cells.Item(2, 9).Value = 2287.0715  # I2: 2323.7693 -> 2287.0715
$ws.Cells.Item(2, 10).Value = 2499  # J2: 2499.5 -> 2499
$ws.Cells.Item(2, 11).Value = 2287.0715  # K2: 2323.7693 -> 2287.0715
$ws.Cells.Item(2, 12).Value = 2499  # L2: 2499.5 -> 2499
$ws.Cells.Item(2, 13).Value = -2174.0715  # M2: -2210.7693 -> -2174.0715
$ws.Cells.Item(2, 14).Value = -2725  # N2: -2725.5 -> -2725
$ws.Cells.Item(74, 8).Value = 2676.2258  # H74: 2729.1292 -> 2676.2258
$ws.Cells.Item(74, 9).Value = 2306.7693  # I74: 2369.8462 -> 2306.7693
$ws.Cells.Item(74, 11).Value = 2306.7693  # K74: 2369.8462 -> 2306.7693
$ws.Cells.Item(74, 13).Value = -1432.7693  # M74: -1495.8462 -> -1432.7693
$ws.Cells.Item(77, 8).Value = 2676.2258  # H77: 2729.1292 -> 2676.2258
$ws.Cells.Item(77, 9).Value = 2306.7693  # I77: 2369.8462 -> 2306.7693
$ws.Cells.Item(77, 11).Value = 11533.8465  # K77: 11849.231 -> 11533.8465
$ws.Cells.Item(77, 13).Value = -7165.8465  # M77: -7481.231 -> -7165.8465
$ws.Cells.Item(97, 8).Value = 561.9091  # H97: 582.1667 -> 561.9091
$ws.Cells.Item(97, 9).Value = 538.1  # I97: 582.1667 -> 538.1
$ws.Cells.Item(97, 10).Value = 800  # J97: 0 -> 800
$ws.Cells.Item(97, 11).Value = 538.1  # K97: 582.1667 -> 538.1
$ws.Cells.Item(97, 12).Value = 800  # L97: 0 -> 800
$ws.Cells.Item(97, 13).Value = -42.10000000000002  # M97: -86.16669999999999 -> -42.10000000000002
$ws.Cells.Item(97, 14).Value = -1792  # N97: None -> -1792
$ws.Cells.Item(116, 8).Value = 2301.2  # H116: 2347.2 -> 2301.2
$ws.Cells.Item(116, 9).Value = 2287.0715  # I116: 2323.7693 -> 2287.0715
$ws.Cells.Item(116, 10).Value = 2499  # J116: 2499.5 -> 2499
$ws.Cells.Item(116, 11).Value = 2287.0715  # K116: 2323.7693 -> 2287.0715
$ws.Cells.Item(116, 12).Value = 2499  # L116: 2499.5 -> 2499
$ws.Cells.Item(116, 13).Value = 6.928499999999985  # M116: -29.76929999999993 -> 6.928499999999985
$ws.Cells.Item(116, 14).Value = -7087  # N116: -7087.5 -> -7087
$ws.Cells.Item(132, 8).Value = 7201.07  # H132: 7401.1323 -> 7201.07
$ws.Cells.Item(132, 9).Value = 5187.987  # I132: 5301.189 -> 5187.987
$ws.Cells.Item(132, 10).Value = 22500.5  # J132: 24667.334 -> 22500.5
$ws.Cells.Item(132, 11).Value = 15563.961  # K132: 15903.567 -> 15563.961
$ws.Cells.Item(132, 12).Value = 67501.5  # L132: 74002.00199999999 -> 67501.5
$ws.Cells.Item(132, 13).Value = -13033.961  # M132: -13373.567 -> -13033.961
$ws.Cells.Item(132, 14).Value = -72561.5  # N132: -79062.00199999999 -> -72561.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2301.2  # H3: 2347.2 -> 2301.2
$ws.Cells.Item(3, 9).Value = 2287.0715  # I3: 2323.7693 -> 2287.0715
$ws.Cells.Item(3, 10).Value = 2499  # J3: 2499.5 -> 2499
$ws.Cells.Item(3, 11).Value = 2287.0715  # K3: 2323.7693 -> 2287.0715
$ws.Cells.Item(3, 12).Value = 2499  # L3: 2499.5 -> 2499
$ws.Cells.Item(3, 13).Value = -2173.0715  # M3: -2209.7693 -> -2173.0715
$ws.Cells.Item(3, 14).Value = -2727  # N3: -2727.5 -> -2727
$ws.Cells.Item(82, 8).Value = 15333.167  # H82: 28571.285 -> 15333.167
$ws.Cells.Item(82, 9).Value = 15333.167  # I82: 15599.8 -> 15333.167
$ws.Cells.Item(82, 10).Value = 0  # J82: 61000 -> 0
$ws.Cells.Item(82, 11).Value = 15333.167  # K82: 15599.8 -> 15333.167
$ws.Cells.Item(82, 12).Value = 0  # L82: 61000 -> 0
$ws.Cells.Item(82, 13).Value = -14950.167  # M82: -15216.8 -> -14950.167
$ws.Cells.Item(82, 14).ClearContents()  # N82: -61766 -> (removed)
$ws.Cells.Item(85, 8).Value = 15333.167  # H85: 28571.285 -> 15333.167
$ws.Cells.Item(85, 9).Value = 15333.167  # I85: 15599.8 -> 15333.167
$ws.Cells.Item(85, 10).Value = 0  # J85: 61000 -> 0
$ws.Cells.Item(85, 11).Value = 15333.167  # K85: 15599.8 -> 15333.167
$ws.Cells.Item(85, 12).Value = 0  # L85: 61000 -> 0
$ws.Cells.Item(85, 13).Value = -14007.167  # M85: -14273.8 -> -14007.167
$ws.Cells.Item(85, 14).ClearContents()  # N85: -63652 -> (removed)
$ws.Cells.Item(94, 8).Value = 45457068  # H94: 66669844 -> 45457068
$ws.Cells.Item(94, 9).Value = 2018.9375  # I94: 2592.4 -> 2018.9375
$ws.Cells.Item(94, 10).Value = 166670530  # J94: 200004340 -> 166670530
$ws.Cells.Item(94, 11).Value = 2018.9375  # K94: 2592.4 -> 2018.9375
$ws.Cells.Item(94, 12).Value = 166670530  # L94: 200004340 -> 166670530
$ws.Cells.Item(94, 13).Value = -1567.9375  # M94: -2141.4 -> -1567.9375
$ws.Cells.Item(94, 14).Value = -166671432  # N94: -200005242 -> -166671432
$ws.Cells.Item(99, 8).Value = 1183  # H99: 83334390 -> 1183
$ws.Cells.Item(99, 9).Value = 1121.9  # I99: 90910110 -> 1121.9
$ws.Cells.Item(99, 10).Value = 1488.5  # J99: 1489 -> 1488.5
$ws.Cells.Item(99, 11).Value = 1121.9  # K99: 90910110 -> 1121.9
$ws.Cells.Item(99, 12).Value = 1488.5  # L99: 1489 -> 1488.5
$ws.Cells.Item(99, 13).Value = 376.0999999999999  # M99: -90908612 -> 376.0999999999999
$ws.Cells.Item(99, 14).Value = -4484.5  # N99: -4485 -> -4484.5
$ws.Cells.Item(105, 8).Value = 1630  # H105: 2218.6667 -> 1630
$ws.Cells.Item(105, 9).Value = 1630  # I105: 2218.6667 -> 1630
$ws.Cells.Item(105, 11).Value = 1630  # K105: 2218.6667 -> 1630
$ws.Cells.Item(105, 13).Value = 117  # M105: -471.6667000000002 -> 117
$ws.Cells.Item(107, 8).Value = 3497.3572  # H107: 4720.154 -> 3497.3572
$ws.Cells.Item(107, 9).Value = 1830.5  # I107: 2086.7 -> 1830.5
$ws.Cells.Item(107, 10).Value = 13498.5  # J107: 13498.333 -> 13498.5
$ws.Cells.Item(107, 11).Value = 1830.5  # K107: 2086.7 -> 1830.5
$ws.Cells.Item(107, 12).Value = 13498.5  # L107: 13498.333 -> 13498.5
$ws.Cells.Item(107, 13).Value = 89.5  # M107: -166.6999999999998 -> 89.5
$ws.Cells.Item(107, 14).Value = -17338.5  # N107: -17338.333 -> -17338.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 47622440  # H86: 52634940 -> 47622440
$ws.Cells.Item(86, 9).Value = 62503064  # I86: 62503130 -> 62503064
$ws.Cells.Item(86, 10).Value = 4445  # J86: 4609.3335 -> 4445
$ws.Cells.Item(86, 11).Value = 62503064  # K86: 62503130 -> 62503064
$ws.Cells.Item(86, 12).Value = 4445  # L86: 4609.3335 -> 4445
$ws.Cells.Item(86, 13).Value = -62501941  # M86: -62502007 -> -62501941
$ws.Cells.Item(86, 14).Value = -6691  # N86: -6855.3335 -> -6691
$ws.Cells.Item(89, 8).Value = 47622440  # H89: 52634940 -> 47622440
$ws.Cells.Item(89, 9).Value = 62503064  # I89: 62503130 -> 62503064
$ws.Cells.Item(89, 10).Value = 4445  # J89: 4609.3335 -> 4445
$ws.Cells.Item(89, 11).Value = 312515320  # K89: 312515650 -> 312515320
$ws.Cells.Item(89, 12).Value = 22225  # L89: 23046.6675 -> 22225
$ws.Cells.Item(89, 13).Value = -312509704  # M89: -312510034 -> -312509704
$ws.Cells.Item(89, 14).Value = -33457  # N89: -34278.6675 -> -33457
$ws.Cells.Item(94, 8).Value = 1416.3846  # H94: 1442.8334 -> 1416.3846
$ws.Cells.Item(94, 10).Value = 1162.8334  # J94: 1175.6 -> 1162.8334
$ws.Cells.Item(94, 12).Value = 1162.8334  # L94: 1175.6 -> 1162.8334
$ws.Cells.Item(94, 14).Value = -2064.8334  # N94: -2077.6 -> -2064.8334
$ws.Cells.Item(132, 8).Value = 1280.7059  # H132: 1455.9412 -> 1280.7059
$ws.Cells.Item(132, 9).Value = 1214.4166  # I132: 1462.6666 -> 1214.4166
$ws.Cells.Item(132, 11).Value = 3643.2498  # K132: 4387.9998 -> 3643.2498
$ws.Cells.Item(132, 13).Value = -1113.2498  # M132: -1857.9998 -> -1113.2498

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 722  # H5: 469.72726 -> 722
$ws.Cells.Item(5, 9).Value = 494.3  # I5: 469.72726 -> 494.3
$ws.Cells.Item(5, 10).Value = 2999  # J5: 0 -> 2999
$ws.Cells.Item(5, 11).Value = 1482.9  # K5: 1409.18178 -> 1482.9
$ws.Cells.Item(5, 12).Value = 8997  # L5: 0 -> 8997
$ws.Cells.Item(5, 13).Value = -1370.9  # M5: -1297.18178 -> -1370.9
$ws.Cells.Item(5, 14).Value = -9221  # N5: None -> -9221
$ws.Cells.Item(8, 8).Value = 2262.6667  # H8: 1552.6 -> 2262.6667
$ws.Cells.Item(8, 9).Value = 2262.6667  # I8: 1552.6 -> 2262.6667
$ws.Cells.Item(8, 11).Value = 6788.000100000001  # K8: 4657.799999999999 -> 6788.000100000001
$ws.Cells.Item(8, 13).Value = -6649.000100000001  # M8: -4518.799999999999 -> -6649.000100000001
$ws.Cells.Item(63, 8).Value = 7482.3335  # H63: 9976 -> 7482.3335
$ws.Cells.Item(63, 9).Value = 6973.5  # I63: 8398.429 -> 6973.5
$ws.Cells.Item(63, 10).Value = 8500  # J63: 15497.5 -> 8500
$ws.Cells.Item(63, 11).Value = 20920.5  # K63: 25195.287 -> 20920.5
$ws.Cells.Item(63, 12).Value = 25500  # L63: 46492.5 -> 25500
$ws.Cells.Item(63, 13).Value = -20171.5  # M63: -24446.287 -> -20171.5
$ws.Cells.Item(63, 14).Value = -26998  # N63: -47990.5 -> -26998
$ws.Cells.Item(66, 8).Value = 7482.3335  # H66: 9976 -> 7482.3335
$ws.Cells.Item(66, 9).Value = 6973.5  # I66: 8398.429 -> 6973.5
$ws.Cells.Item(66, 10).Value = 8500  # J66: 15497.5 -> 8500
$ws.Cells.Item(66, 11).Value = 62761.5  # K66: 75585.861 -> 62761.5
$ws.Cells.Item(66, 12).Value = 76500  # L66: 139477.5 -> 76500
$ws.Cells.Item(66, 13).Value = -59017.5  # M66: -71841.861 -> -59017.5
$ws.Cells.Item(66, 14).Value = -83988  # N66: -146965.5 -> -83988
$ws.Cells.Item(74, 8).Value = 19750  # H74: 18291.666 -> 19750
$ws.Cells.Item(74, 10).Value = 19687.5  # J74: 17950 -> 19687.5
$ws.Cells.Item(74, 12).Value = 59062.5  # L74: 53850 -> 59062.5
$ws.Cells.Item(74, 14).Value = -61184.5  # N74: -55972 -> -61184.5
$ws.Cells.Item(77, 8).Value = 19750  # H77: 18291.666 -> 19750
$ws.Cells.Item(77, 10).Value = 19687.5  # J77: 17950 -> 19687.5
$ws.Cells.Item(77, 12).Value = 177187.5  # L77: 161550 -> 177187.5
$ws.Cells.Item(77, 14).Value = -187795.5  # N77: -172158 -> -187795.5
$ws.Cells.Item(131, 8).Value = 1455.6364  # H131: 1598.6 -> 1455.6364
$ws.Cells.Item(131, 9).Value = 915.4  # I131: 940.4 -> 915.4
$ws.Cells.Item(131, 10).Value = 1905.8334  # J131: 2915 -> 1905.8334
$ws.Cells.Item(131, 11).Value = 2746.2  # K131: 2821.2 -> 2746.2
$ws.Cells.Item(131, 12).Value = 5717.5002  # L131: 8745 -> 5717.5002
$ws.Cells.Item(131, 13).Value = 2293.8  # M131: 2218.8 -> 2293.8
$ws.Cells.Item(131, 14).Value = -15797.5002  # N131: -18825 -> -15797.5002
$ws.Cells.Item(135, 8).Value = 722  # H135: 469.72726 -> 722
$ws.Cells.Item(135, 9).Value = 494.3  # I135: 469.72726 -> 494.3
$ws.Cells.Item(135, 10).Value = 2999  # J135: 0 -> 2999
$ws.Cells.Item(135, 11).Value = 4448.7  # K135: 4227.54534 -> 4448.7
$ws.Cells.Item(135, 12).Value = 26991  # L135: 0 -> 26991
$ws.Cells.Item(135, 13).Value = -1913.7  # M135: -1692.54534 -> -1913.7
$ws.Cells.Item(135, 14).Value = -32061  # N135: None -> -32061

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3892.88  # H80: 3912.074 -> 3892.88
$ws.Cells.Item(80, 9).Value = 3684.3333  # I80: 3851 -> 3684.3333
$ws.Cells.Item(80, 10).Value = 4085.3845  # J80: 3960.9333 -> 4085.3845
$ws.Cells.Item(80, 11).Value = 3684.3333  # K80: 3851 -> 3684.3333
$ws.Cells.Item(80, 12).Value = 4085.3845  # L80: 3960.9333 -> 4085.3845
$ws.Cells.Item(80, 13).Value = -2686.3333  # M80: -2853 -> -2686.3333
$ws.Cells.Item(80, 14).Value = -6081.3845  # N80: -5956.933300000001 -> -6081.3845
$ws.Cells.Item(83, 8).Value = 3892.88  # H83: 3912.074 -> 3892.88
$ws.Cells.Item(83, 9).Value = 3684.3333  # I83: 3851 -> 3684.3333
$ws.Cells.Item(83, 10).Value = 4085.3845  # J83: 3960.9333 -> 4085.3845
$ws.Cells.Item(83, 11).Value = 18421.6665  # K83: 19255 -> 18421.6665
$ws.Cells.Item(83, 12).Value = 20426.9225  # L83: 19804.6665 -> 20426.9225
$ws.Cells.Item(83, 13).Value = -13429.6665  # M83: -14263 -> -13429.6665
$ws.Cells.Item(83, 14).Value = -30410.9225  # N83: -29788.6665 -> -30410.9225
$ws.Cells.Item(107, 8).Value = 1623.1538  # H107: 1629.8462 -> 1623.1538
$ws.Cells.Item(107, 9).Value = 1374.6364  # I107: 1382.5454 -> 1374.6364
$ws.Cells.Item(107, 11).Value = 1374.6364  # K107: 1382.5454 -> 1374.6364
$ws.Cells.Item(107, 13).Value = 545.3635999999999  # M107: 537.4546 -> 545.3635999999999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2838.2903  # H16: 2891.3333 -> 2838.2903
$ws.Cells.Item(16, 9).Value = 2792.7273  # I16: 2866.3333 -> 2792.7273
$ws.Cells.Item(16, 11).Value = 2792.7273  # K16: 2866.3333 -> 2792.7273
$ws.Cells.Item(16, 13).Value = -2622.7273  # M16: -2696.3333 -> -2622.7273
$ws.Cells.Item(61, 8).Value = 948.9167  # H61: 976.5454999999999 -> 948.9167
$ws.Cells.Item(61, 9).Value = 833.1818  # I61: 852 -> 833.1818
$ws.Cells.Item(61, 11).Value = 833.1818  # K61: 852 -> 833.1818
$ws.Cells.Item(61, 13).Value = -631.1818  # M61: -650 -> -631.1818
$ws.Cells.Item(68, 8).Value = 3820.9285  # H68: 3999.5386 -> 3820.9285
$ws.Cells.Item(68, 9).Value = 2322.818  # I68: 2405.2 -> 2322.818
$ws.Cells.Item(68, 11).Value = 2322.818  # K68: 2405.2 -> 2322.818
$ws.Cells.Item(68, 13).Value = -1573.818  # M68: -1656.2 -> -1573.818
$ws.Cells.Item(71, 8).Value = 3820.9285  # H71: 3999.5386 -> 3820.9285
$ws.Cells.Item(71, 9).Value = 2322.818  # I71: 2405.2 -> 2322.818
$ws.Cells.Item(71, 11).Value = 11614.09  # K71: 12026 -> 11614.09
$ws.Cells.Item(71, 13).Value = -7870.09  # M71: -8282 -> -7870.09
$ws.Cells.Item(93, 8).Value = 4417.273  # H93: 4529.7 -> 4417.273
$ws.Cells.Item(93, 9).Value = 3510  # I93: 3537.125 -> 3510
$ws.Cells.Item(93, 11).Value = 3510  # K93: 3537.125 -> 3510
$ws.Cells.Item(93, 13).Value = -2262  # M93: -2289.125 -> -2262
$ws.Cells.Item(113, 8).Value = 948.9167  # H113: 976.5454999999999 -> 948.9167
$ws.Cells.Item(113, 9).Value = 833.1818  # I113: 852 -> 833.1818
$ws.Cells.Item(113, 11).Value = 833.1818  # K113: 852 -> 833.1818
$ws.Cells.Item(113, 13).Value = 1336.8182  # M113: 1318 -> 1336.8182

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 6174694.5  # H96: 7409214.5 -> 6174694.5
$ws.Cells.Item(96, 9).Value = 9260963  # I96: 12347253 -> 9260963
$ws.Cells.Item(96, 11).Value = 9260963  # K96: 12347253 -> 9260963
$ws.Cells.Item(96, 13).Value = -9259590  # M96: -12345880 -> -9259590
